$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

function Set-PlainValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

Set-PlainValue "D2" "27.473.30"
Set-PlainValue "E2" "  +1.58%  "
Set-PlainValue "D3" "1.565.04"
Set-PlainValue "E3" "  -0.04%  "
Set-TextValue "D4" "0.992"
Set-PlainValue "E4" "  -1.24%  "
Set-TextValue "D5" "210.72"
Set-PlainValue "E5" "  +1.01%  "
Set-PlainValue "E6" "  -0.48%  "
Set-TextValue "D7" "0.991"
Set-PlainValue "E7" "  -1.35%  "
Set-TextValue "D8" "22.58"
Set-PlainValue "E8" "  +2.05%  "
Set-PlainValue "E9" "  +0.08%  "
Set-TextValue "D10" "0.0594"
Set-PlainValue "E10" "  -0.63%  "
Set-TextValue "D11" "0.0870"
Set-PlainValue "E11" "  +1.34%  "
Set-PlainValue "D12" "1.789.96"
Set-PlainValue "E12" "  +0.08%  "
Set-PlainValue "D13" "1.558.22"
Set-PlainValue "E13" "  -0.42%  "
Set-TextValue "D14" "3.76"
Set-PlainValue "E14" "  -0.33%  "
Set-PlainValue "E15" "  -0.44%  "
Set-PlainValue "D16" "27.500.23"
Set-PlainValue "E16" "  +1.76%  "
Set-TextValue "D17" "62.34"
Set-PlainValue "E17" "  +0.71%  "
Set-TextValue "D18" "224.87"
Set-PlainValue "E18" "  +4.08%  "
Set-TextValue "D19" "7.49"
Set-PlainValue "E19" "  +0.94%  "
Set-PlainValue "D20" "0.0₃0703"
Set-PlainValue "E20" "  -0.45%  "
Set-TextValue "D21" "0.991"
Set-PlainValue "E21" "  -1.41%  "
Set-TextValue "D22" "4.11"
Set-PlainValue "E22" "  -0.98%  "
Set-TextValue "D23" "9.41"
Set-PlainValue "E23" "  +2.01%  "
Set-PlainValue "E24" "  +0.30%  "
Set-TextValue "D25" "149.75"
Set-PlainValue "E25" "  -2.88%  "
Set-TextValue "D26" "15.13"
Set-PlainValue "E26" "  +0.35%  "
Set-TextValue "D27" "6.60"
Set-PlainValue "E27" "  -0.35%  "
Set-PlainValue "E28" "  +1.54%  "
Set-TextValue "D29" "0.992"
Set-PlainValue "E29" "  -1.30%  "
Set-PlainValue "E30" "  +1.09%  "
Set-PlainValue "E31" "  -0.83%  "
Set-TextValue "D32" "3.23"
Set-PlainValue "E32" "  -0.02%  "
Set-PlainValue "D33" "1.445.21"
Set-PlainValue "E33" "  +1.14%  "
Set-TextValue "D34" "3.14"
Set-PlainValue "E34" "  -2.45%  "
Set-PlainValue "E35" "  +2.15%  "
Set-PlainValue "E36" "  -0.73%  "
Set-PlainValue "E37" "  -0.83%  "
Set-TextValue "D38" "0.0167"
Set-PlainValue "E38" "  -0.25%  "
Set-TextValue "D39" "0.540"
Set-PlainValue "E39" "  +1.03%  "
Set-TextValue "D40" "0.812"
Set-PlainValue "E40" "  -0.37%  "
Set-PlainValue "B41" "FraxShare"
Set-PlainValue "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D41" "5.73"
Set-PlainValue "E41" "  -1.20%  "
Set-PlainValue "B42" "MXToken"
Set-PlainValue "C42" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D42" "2.36"
Set-PlainValue "E42" "  +1.63%  "
Set-TextValue "D43" "0.991"
Set-PlainValue "E43" "  -1.45%  "
Set-TextValue "D44" "1.83"
Set-PlainValue "E44" "  +5.03%  "
Set-TextValue "D45" "0.975"
Set-PlainValue "E45" "  -2.97%  "
Set-TextValue "D46" "64.38"
Set-PlainValue "E46" "  -0.73%  "
Set-PlainValue "D47" "1.702.92"
Set-TextValue "D48" "86.61"
Set-PlainValue "E48" "  -0.10%  "
Set-TextValue "D49" "0.0525"
Set-PlainValue "E49" "  +1.14%  "
Set-PlainValue "B50" "BabyDogeCoin"
Set-PlainValue "C50" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-PlainValue "D50" "0.0₇0960"
Set-PlainValue "E50" "  -6.30%  "
Set-PlainValue "B51" "Algorand"
Set-PlainValue "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.0945"
Set-PlainValue "E51" "  -1.69%  "
